# Reduce the page margins on the document's (only) section from the Normal
# default of 1 inch (1440 twips / 72 pt) to 0.5 inch (720 twips / 36 pt) on
# each of the top/right/bottom/left edges. Header/footer distances and the
# gutter are left untouched, matching the target <w:pgMar/> values of
# top="720" right="720" bottom="720" left="720" header="720" footer="720"
# gutter="0".
$d = $word.ActiveDocument

foreach ($section in $d.Sections) {
    $pageSetup = $section.PageSetup
    $pageSetup.TopMargin = 36
    $pageSetup.BottomMargin = 36
    $pageSetup.LeftMargin = 36
    $pageSetup.RightMargin = 36
}
